$d = $word.ActiveDocument
$wdParagraph = 4

# ---------------------------------------------------------------------
# 1) "Fridlysta arter" section: merge the "Följande fridlysta arter..."
#    paragraph with the "Revlummer (§9)" bullet paragraph into a single
#    paragraph, lower-casing "revlummer" and adding a trailing period.
# ---------------------------------------------------------------------

# Remove the whole "Revlummer (§9)" bullet paragraph (match case so we
# target the capitalised heading-style occurrence, not the lower-case
# mention earlier in the document).
$bulletRng = $d.Content
$bulletRng.Find.Execute("Revlummer (§9)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bulletRng.Expand($wdParagraph) | Out-Null
$bulletRng.Delete()

# Update the remaining intro paragraph's text.
$introRng = $d.Content
$introRng.Find.Execute("avverkningsanmälda skogen: ", $true, $false, $false, $false, $false, $true, 1, $false, "avverkningsanmälda skogen: revlummer (§9).", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Remove italic formatting from the "I den avverkningsanmälda
#    skogen..." comment run (first occurrence only) and append a period.
#    Re-insert as plain (unformatted) text so no explicit rPr remains.
# ---------------------------------------------------------------------
$commentText = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"
$found = $d.Content
$found.Find.Execute($commentText) | Out-Null
$start = $found.Start
$end = $found.End
$insertPoint = $d.Range($end, $end)
$insertPoint.InsertAfter($commentText + ".")
$oldRange = $d.Range($start, $end)
$oldRange.Delete()

# ---------------------------------------------------------------------
# 3) Drop the trailing space in the "Certifikatsinnehavaren ska skydda
#    sällsynta arter..." (6.4) run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas. ", $true, $false, $false, $false, $false, $true, 1, $false, "arter utanför skogsbruksenhetens gränser när beslut om åtgärder inom skogsbruksenheten ska fattas.", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Collapse the two "6.4.1" paragraphs + the intervening nyckelbiotoper
#    bullet into a single "6.4.3" paragraph.
# ---------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("Följande biotoper undantas från alla skogsbruksåtgärder") | Out-Null
$f1.Expand($wdParagraph) | Out-Null
$s1 = $f1.Start

$f2 = $d.Content
$f2.Find.Execute("nyckelbiotoper enligt Skogsstyrelsens definition och metod (1995)") | Out-Null
$f2.Expand($wdParagraph) | Out-Null
$e2 = $f2.End

$unionRng = $d.Range($s1, $e2)
$unionRng.Delete()

$d.Content.Find.Execute("6.4.1 ", $true, $false, $false, $false, $false, $true, 1, $false, "6.4.3 ", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Update the date in the first-page header.
# ---------------------------------------------------------------------
$sec = $d.Sections.First
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists) {
        $h.Range.Find.Execute("2023-10-22", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-25", 2) | Out-Null
    }
}
